# Generate Report for handoff
# Updates the "b.md.md" row on each sheet to reflect that a new handoff
# package (b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0) is now ready.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet: row 3 is the "b.md.md" row. Its zh-cn/de-de status
# moves from "Handed back: in sync with en-US" to "Ready for handoff".
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------
# zh-cn sheet: row 3 ("b.md.md") gets a new handoff file/time, and its
# status is updated too. The hyperlink text shown for the handoff file
# cell needs to be updated to match (its target URL is left untouched).
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-01-26 05:35:07"
foreach ($hl in $zhcn.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------
# de-de sheet: same kind of update as zh-cn, with the de-de handoff file.
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
$dede.Range("D3").Value = "2016-01-26 05:35:17"
foreach ($hl in $dede.Hyperlinks) {
    if ($hl.Range.Address() -eq '$C$3') {
        $hl.TextToDisplay = "b.md.b3a40d6229ff1a8b48804fcfc66c95922eb78fd0.de-de.xlf"
    }
}
